$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (existing H..N shift right to I..O)
$ws.Columns.Item(8).Insert()

# New header for the inserted column
$ws.Range("H1").Value = "CO2/(CO+CO2)"

# New ratio formula: CO2/(CO+CO2) using feed-derived columns E (Yco) and F (Yco2)
$ws.Range("H2").Formula = "=F2/(E2+F2)"
$ws.Range("H3:H28").Formula = "=F3/(E3+F3)"

# Match the formatting of the new column to the rest of the header/data style
$ws.Range("H1").Copy()
$ws.Range("H2:H28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the new column selected, matching the edited workbook's view state
$ws.Range("H2:H28").Select()
